# feat: add 2022-Q3 data
#
# - Duplicate the existing "2022-Q2" sheet (placed right after it) so the
#   original Q2 fund-holding data is preserved on its own tab.
# - Rename the original sheet to "2022-Q3" and overwrite its data with the
#   new quarter's fund-holding figures (re-using the "总计" sheet's header /
#   index-column formatting).
# - Rename the duplicate back to "2022-Q2".
# - Update the "总计" (summary) sheet: the existing row now reflects Q3, and
#   a new row is appended below it with the old Q2 figures.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Duplicate "2022-Q2" -> becomes the new trailing "2022-Q2" tab, then
#    repurpose the original tab (same sheetId/rId) as "2022-Q3".
# ---------------------------------------------------------------------
$wsQ2 = $wb.Worksheets.Item("2022-Q2")
$wsQ2.Copy($null, $wsQ2)

$wsQ2.Name = "2022-Q3"

$wsCopy = $wb.Worksheets.Item("2022-Q2 (2)")
$wsCopy.Name = "2022-Q2"

# ---------------------------------------------------------------------
# 2) Overwrite the (now named) "2022-Q3" sheet with the new quarter data.
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsTotal = $wb.Worksheets.Item("总计")

$fundRows = @(
    @(0, "014273", "广发北交所精选两年定开混合A", "3.37", "64.25", "7.50", "0.2528", 3),
    @(1, "014277", "万家北交所慧选两年定期开放混合A", "3.56", "93.97", "6.89", "0.2453", 3),
    @(2, "014271", "大成北交所两年定开混合A", "3.45", "65.31", "6.84", "0.2360", 5),
    @(3, "014279", "汇添富北交所创新精选两年定开混合A", "3.20", "93.27", "4.95", "0.1584", 3),
    @(4, "014294", "南方北交所精选两年定开混合", "4.26", "75.23", "2.91", "0.1240", 7),
    @(5, "014274", "广发北交所精选两年定开混合C", "0.85", "64.25", "7.50", "0.0638", 3),
    @(6, "014272", "大成北交所两年定开混合C", "0.82", "65.31", "6.84", "0.0561", 5),
    @(7, "014278", "万家北交所慧选两年定期开放混合C", "0.49", "93.97", "6.89", "0.0338", 3),
    @(8, "014280", "汇添富北交所创新精选两年定开混合C", "0.51", "93.27", "4.95", "0.0252", 3)
)

foreach ($fr in $fundRows) {
    $r = $fr[0] + 2
    $wsQ3.Cells.Item($r, 1).Value = $fr[0]
    $wsQ3.Cells.Item($r, 2).Value = "'" + $fr[1]
    $wsQ3.Cells.Item($r, 3).Value = $fr[2]
    $wsQ3.Cells.Item($r, 4).Value = "'" + $fr[3]
    $wsQ3.Cells.Item($r, 5).Value = "'" + $fr[4]
    $wsQ3.Cells.Item($r, 6).Value = "'" + $fr[5]
    $wsQ3.Cells.Item($r, 7).Value = "'" + $fr[6]
    $wsQ3.Cells.Item($r, 8).Value = $fr[7]
}

# The new quarter's sheet reuses the "总计" sheet's header / index-column
# look (bold, centered, thin border) rather than the old sheet's style.
$wsTotal.Range("B1").Copy()
$wsQ3.Range("B1:H1").PasteSpecial(-4122)

$wsTotal.Range("A2").Copy()
$wsQ3.Range("A2:A10").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 3) Update the "总计" summary sheet: row 2 -> Q3 figures, new row 3 ->
#    old Q2 figures (copy formatting from row 2's A cell).
# ---------------------------------------------------------------------
$wsTotal.Range("A2").Copy()
$wsTotal.Range("A3").PasteSpecial(-4122)

$wsTotal.Cells.Item(3, 1).Value = 1
$wsTotal.Cells.Item(3, 2).Value = "2022-Q2"
$wsTotal.Cells.Item(3, 3).Value = 9
$wsTotal.Cells.Item(3, 4).Value = 1.46

$wsTotal.Cells.Item(2, 2).Value = "2022-Q3"
$wsTotal.Cells.Item(2, 3).Value = 9
$wsTotal.Cells.Item(2, 4).Value = 1.2

$excel.CutCopyMode = $false

# Keep "总计" as the active sheet (unchanged from the source workbook).
$wsTotal.Activate()
